# "adding averages and more checks"
#
# xl/styles.xml is shared across the whole workbook, so the font-table
# edit (the old 14pt "title" font and the plain bold "header" font are
# consolidated into a single bold + white font, used by both the report
# title row and the header row) shows up on every sheet that uses those
# two cell styles - not just "Training Dashboard".
#
#   1. Title row (row 1) and header row (row 2) on every sheet: font
#      becomes bold white (instead of bold/black, and the title loses its
#      old 14pt size).
#   2. Training Dashboard!H3 ("PERIOD TO EXPIRE"): -43 -> -51
#   3. Training Dashboard!I3 ("LAST UPDATE"): "08-Sep-2025" -> "16-Sep-2025"
#      (kept as literal text, not reinterpreted as a date serial)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastCol = $used.Columns.Count

    # Row 1 = report title (only the merged-region anchor cell, A1, is an
    # actual cell in the sheet). Color is set first so no intermediate
    # state happens to land on a pre-existing font/style combination.
    $title = $ws.Range("A1")
    $title.Font.Color = 16777215
    $title.Font.Bold = $true
    $title.Font.Size = 11

    # Row 2 = column headers, spanning every used column.
    $header = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $header.Font.Color = 16777215
    $header.Font.Bold = $true
}

# --- Data updates on the Training Dashboard sheet, row 3 ---
$ws = $wb.Worksheets.Item("Training Dashboard")

$hCell = $ws.Range("H3")
$hCell.Value = -51

$iCell = $ws.Range("I3")
$iCell.Formula = '="16-Sep-2025"'
$iCell.Copy()
$iCell.PasteSpecial(-4163)
